$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1672
$ws.Range("I6").Value = 241
$ws.Range("K6").Value = 723
$ws.Range("M6").Value = -611
$ws.Range("N6").ClearContents()
$ws.Range("H9").Value = 1983.1428
$ws.Range("I9").Value = 308.22223
$ws.Range("J9").Value = 3239.3333
$ws.Range("K9").Value = 308.22223
$ws.Range("L9").Value = 3239.3333
$ws.Range("M9").Value = -139.22223
$ws.Range("N9").Value = -3577.3333
$ws.Range("H12").Value = 1054.9
$ws.Range("I12").Value = 3000
$ws.Range("J12").Value = 838.7778
$ws.Range("K12").Value = 3000
$ws.Range("L12").Value = 838.7778
$ws.Range("M12").Value = -2830
$ws.Range("N12").Value = -1178.7778
$ws.Range("H94").Value = 6350.8
$ws.Range("I94").Value = 6350.8
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 6350.8
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -5899.8
$ws.Range("N94").ClearContents()
$ws.Range("H113").Value = 2280
$ws.Range("I113").Value = 2400
$ws.Range("K113").Value = 2400
$ws.Range("M113").Value = 854
$ws.Range("N113").ClearContents()
$ws.Range("H137").Value = 3116.889
$ws.Range("I137").Value = 776.5
$ws.Range("K137").Value = 2329.5
$ws.Range("M137").Value = 220.5
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3463.2917
$ws.Range("I32").Value = 2309.6956
$ws.Range("K32").Value = 2309.6956
$ws.Range("M32").Value = -2022.6956
$ws.Range("N32").ClearContents()
$ws.Range("H61").Value = 2500
$ws.Range("I61").Value = 2500
$ws.Range("K61").Value = 2500
$ws.Range("M61").Value = -2288
$ws.Range("H115").Value = 20728
$ws.Range("I115").Value = 10000
$ws.Range("J115").Value = 26092
$ws.Range("K115").Value = 10000
$ws.Range("L115").Value = 26092
$ws.Range("M115").Value = -8433
$ws.Range("N115").Value = -29226
$ws.Range("H132").Value = 3242.75
$ws.Range("I132").Value = 2212.375
$ws.Range("J132").Value = 5303.5
$ws.Range("K132").Value = 6637.125
$ws.Range("L132").Value = 15910.5
$ws.Range("M132").Value = -4107.125
$ws.Range("N132").Value = -20970.5
$ws.Range("H136").Value = 2500
$ws.Range("I136").Value = 2500
$ws.Range("K136").Value = 7500
$ws.Range("M136").Value = -4950

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2410.7144
$ws.Range("I99").Value = 1995.8334
$ws.Range("K99").Value = 1995.8334
$ws.Range("M99").Value = -497.8334
$ws.Range("N99").ClearContents()
$ws.Range("H105").Value = 4592.3335
$ws.Range("I105").Value = 2831.3333
$ws.Range("J105").Value = 5766.3335
$ws.Range("K105").Value = 2831.3333
$ws.Range("L105").Value = 5766.3335
$ws.Range("M105").Value = -1084.3333
$ws.Range("N105").Value = -9260.333500000001
$ws.Range("H114").Value = 40000
$ws.Range("J114").Value = 40000
$ws.Range("L114").Value = 40000
$ws.Range("N114").Value = -48678
$ws.Range("H134").Value = 2348.75
$ws.Range("I134").Value = 2093.5557
$ws.Range("K134").Value = 6280.6671
$ws.Range("M134").Value = -3745.6671
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H18").Value = 79965.55
$ws.Range("J18").Value = 79965.55
$ws.Range("L18").Value = 79965.55
$ws.Range("N18").Value = -80425.55
$ws.Range("H114").Value = 105995
$ws.Range("J114").Value = 105995
$ws.Range("L114").Value = 105995
$ws.Range("N114").Value = -114673
$ws.Range("H132").Value = 8694.643
$ws.Range("J132").Value = 11869.777
$ws.Range("L132").Value = 35609.331
$ws.Range("N132").Value = -40669.331
$ws.Range("H134").Value = 2812
$ws.Range("I134").Value = 2546.4167
$ws.Range("J134").Value = 5999
$ws.Range("K134").Value = 7639.250100000001
$ws.Range("L134").Value = 17997
$ws.Range("M134").Value = -5104.250100000001
$ws.Range("N134").Value = -23067
$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("M140").ClearContents()
$ws.Range("H141").Value = 442210
$ws.Range("J141").Value = 442210
$ws.Range("L141").Value = 442210
$ws.Range("N141").Value = -452570

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 740
$ws.Range("I7").Value = 100
$ws.Range("K7").Value = 300
$ws.Range("M7").Value = -188
$ws.Range("N7").ClearContents()
$ws.Range("H12").Value = 177087.3
$ws.Range("I12").Value = 428605.28
$ws.Range("K12").Value = 1285815.84
$ws.Range("M12").Value = -1285642.84
$ws.Range("N12").ClearContents()
$ws.Range("H33").Value = 199.75
$ws.Range("I33").Value = 200
$ws.Range("K33").Value = 1200
$ws.Range("M33").Value = -917
$ws.Range("N33").ClearContents()
$ws.Range("H64").Value = 500
$ws.Range("I64").Value = 500
$ws.Range("K64").Value = 1500
$ws.Range("M64").Value = -1230
$ws.Range("H67").Value = 500
$ws.Range("I67").Value = 500
$ws.Range("K67").Value = 1500
$ws.Range("M67").Value = -564
$ws.Range("H68").Value = 999
$ws.Range("J68").Value = 999
$ws.Range("L68").Value = 2997
$ws.Range("N68").Value = -4619
$ws.Range("H71").Value = 999
$ws.Range("J71").Value = 999
$ws.Range("L71").Value = 8991
$ws.Range("N71").Value = -17103
$ws.Range("H80").Value = 4599.2
$ws.Range("I80").Value = 1999.5
$ws.Range("J80").Value = 14998
$ws.Range("K80").Value = 5998.5
$ws.Range("L80").Value = 44994
$ws.Range("M80").Value = -5062.5
$ws.Range("N80").Value = -46866
$ws.Range("H83").Value = 4599.2
$ws.Range("I83").Value = 1999.5
$ws.Range("J83").Value = 14998
$ws.Range("K83").Value = 17995.5
$ws.Range("L83").Value = 134982
$ws.Range("M83").Value = -13315.5
$ws.Range("N83").Value = -144342

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H108").Value = 105995
$ws.Range("J108").Value = 105995
$ws.Range("L108").Value = 105995
$ws.Range("N108").Value = -113675
$ws.Range("H132").Value = 3175.1765
$ws.Range("I132").Value = 2548.9
$ws.Range("J132").Value = 4069.8572
$ws.Range("K132").Value = 7646.700000000001
$ws.Range("L132").Value = 12209.5716
$ws.Range("M132").Value = -5116.700000000001
$ws.Range("N132").Value = -17269.5716

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H55").Value = 903.9545000000001
$ws.Range("J55").Value = 1599.7273
$ws.Range("L55").Value = 1599.7273
$ws.Range("N55").Value = -1945.7273
$ws.Range("H122").Value = 2399.6667
$ws.Range("J122").Value = 2499.5
$ws.Range("L122").Value = 7498.5
$ws.Range("N122").Value = -12398.5
$ws.Range("H136").Value = 4299.75
$ws.Range("I136").Value = 3969.7
$ws.Range("K136").Value = 11909.1
$ws.Range("M136").Value = -9359.099999999999
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H114").Value = 41000
$ws.Range("J114").Value = 41000
$ws.Range("L114").Value = 41000
$ws.Range("N114").Value = -49678
$ws.Range("H132").Value = 4222
$ws.Range("I132").Value = 4222
$ws.Range("K132").Value = 12666
$ws.Range("M132").Value = -10136
